$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: add B24, C24, D24 following the existing pattern used in rows 15-23,
# and extend the E24 formula to include the new B24/C24/D24 concat (shared formula with row 24 itself).
$ws.Range("B24").Value = " --add-data=""assets\"
$ws.Range("C24").Value = "play.png"
$ws.Range("D24").Value = ";assets"""
$ws.Range("B24").Style = $ws.Range("B15").Style

$ws.Range("E24").Formula = "=_xlfn.CONCAT(E23,B24,C24,D24)"

# New row 25: append the final builderBase.py command
$ws.Range("E25").Formula = "=_xlfn.CONCAT(E24,"" builderBase.py"")"

$ws.Range("E25").Select()

$wb.Save()
